# Generate Report for Archive
# - Flip the localization status from "Ready for handoff" to "In Translation"
#   everywhere it is reported (Overview summary columns + per-locale Status
#   columns), then re-fit the now-shorter Status columns to their content
#   (they auto-sized narrower since "In Translation" renders tighter than
#   "Ready for handoff").

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

# Target OOXML column width (<col width="...">) for the shrunk Status
# columns, taken from the regenerated report: 13.4101848602295 characters.
# ColumnWidth assignments get pixel-snapped to that stored width on write,
# so request a value from the middle of the input bucket that snaps to the
# stored width closest to the target (13.3333.. char, i.e. 96px @ 7px/char).
$newWidth = 12.5

# --- Overview sheet: zh-cn (E) / de-de (F) status columns, rows 2-3 ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2:F3").Value = $newStatus
$overview.Range("E1:F1").ColumnWidth = $newWidth

# --- Per-locale detail sheets: Status column (C), rows 2-3 ---
foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("C2:C3").Value = $newStatus
    $ws.Range("C1").ColumnWidth = $newWidth
}
